$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (exhibitions) - insert a brand-new event row at row 32,
# pushing the former rows 32-35 down to 33-36, and bump a handful of
# "want to go" (F column) counters that were already above the insertion
# point.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# Insert a new blank row at position 32 (shifts existing rows 32..35 -> 33..36)
$ws1.Rows.Item(32).Insert()

# Copy formatting from the row that is now directly below (old row 32, now
# at row 33) onto the new blank row so style/borders match the rest of the
# table (matches what Excel does visually when a row is inserted).
$ws1.Cells.Item(33, 1).EntireRow.Copy()
$ws1.Cells.Item(32, 1).PasteSpecial(-4122)

# Fill in the new event's data. Columns B and E hold date-like text that
# Excel would otherwise auto-convert to a real date/time value, so (exactly
# like a human avoiding that auto-conversion in the Excel UI) we lead with
# an apostrophe to force literal text while keeping a plain "General" number
# format.
$ws1.Cells.Item(32, 1).Value = 31
$ws1.Cells.Item(32, 2).Value = "'2024-10-05"
$ws1.Cells.Item(32, 3).Value = "杭州·首届CCPC动漫嘉年华"
$ws1.Cells.Item(32, 4).Value = "长乐路29号五组2幢 杭州运河文化发布中心"
$ws1.Cells.Item(32, 5).Value = "'2024.10.05 09:00-10.06 18:00"
$ws1.Cells.Item(32, 6).Value = 9
$ws1.Cells.Item(32, 7).Value = 39
$ws1.Cells.Item(32, 8).Value = "https://show.bilibili.com/platform/detail.html?id=91102"
$ws1.Cells.Item(32, 9).Value = "//i2.hdslb.com/bfs/openplatform/202408/cf8ib7Q91724210459091.jpeg"

# "Want to go" counter bumps for rows above the insertion point (unaffected
# by the shift).
$ws1.Cells.Item(3, 6).Value = 410
$ws1.Cells.Item(4, 6).Value = 1164
$ws1.Cells.Item(8, 6).Value = 1073
$ws1.Cells.Item(13, 6).Value = 317
$ws1.Cells.Item(14, 6).Value = 360
$ws1.Cells.Item(15, 6).Value = 35
$ws1.Cells.Item(17, 6).Value = 522
$ws1.Cells.Item(18, 6).Value = 1463
$ws1.Cells.Item(19, 6).Value = 5684
$ws1.Cells.Item(21, 6).Value = 1590
$ws1.Cells.Item(22, 6).Value = 378
$ws1.Cells.Item(23, 6).Value = 43
$ws1.Cells.Item(24, 6).Value = 27
$ws1.Cells.Item(25, 6).Value = 5212
$ws1.Cells.Item(26, 6).Value = 5212
$ws1.Cells.Item(27, 6).Value = 127
$ws1.Cells.Item(28, 6).Value = 84
$ws1.Cells.Item(29, 6).Value = 1524

# Counter bumps for the rows that shifted down by one (now at 33..36).
$ws1.Cells.Item(33, 6).Value = 662
$ws1.Cells.Item(34, 6).Value = 86
$ws1.Cells.Item(36, 6).Value = 3808

# ---------------------------------------------------------------------------
# Sheet "演出" (performances)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(4, 6).Value = 20
$ws2.Cells.Item(5, 6).Value = 154
$ws2.Cells.Item(8, 6).Value = 132

# ---------------------------------------------------------------------------
# Sheet "本地生活" (local life)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 9397
$ws3.Cells.Item(4, 6).Value = 2143

# ---------------------------------------------------------------------------
# Sheet "全部类型" (all types) - merged/sorted view of the other three
# sheets; same counters, different row numbers.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value = 9397
$ws4.Cells.Item(4, 6).Value = 2143
$ws4.Cells.Item(6, 6).Value = 410
$ws4.Cells.Item(7, 6).Value = 1164
$ws4.Cells.Item(11, 6).Value = 1073
$ws4.Cells.Item(14, 6).Value = 317
$ws4.Cells.Item(15, 6).Value = 360
$ws4.Cells.Item(16, 6).Value = 35
$ws4.Cells.Item(21, 6).Value = 522
$ws4.Cells.Item(22, 6).Value = 1463
$ws4.Cells.Item(23, 6).Value = 5684
$ws4.Cells.Item(25, 6).Value = 1590
$ws4.Cells.Item(28, 6).Value = 378
$ws4.Cells.Item(31, 6).Value = 5212
$ws4.Cells.Item(32, 6).Value = 5212
$ws4.Cells.Item(33, 6).Value = 127
$ws4.Cells.Item(34, 6).Value = 84
$ws4.Cells.Item(35, 6).Value = 1524
$ws4.Cells.Item(38, 6).Value = 662
$ws4.Cells.Item(39, 6).Value = 86
$ws4.Cells.Item(47, 6).Value = 3808

Write-Output "edit complete"
